$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table cell text: "Wartość sprzedaży minus zużycie"
#    ->  four separate runs (same Arial/bold/12pt formatting) reading
#        "Wartość " + "plus opłat" + "a" + " oraz prowizja"
# ---------------------------------------------------------------------------

$findRange = $d.Content
$found = $findRange.Find.Execute("Wartość sprzedaży minus zużycie")

if ($found) {
    $start = $findRange.Start
    $end   = $findRange.End

    # Segment 1 overwrites the whole original run's text - keeps the
    # original run's rPr (Arial, bold, 12pt) since it rewrites that run in place.
    $seg1 = $d.Range($start, $end)
    $seg1.Text = "Wartość "
    $seg1Start = $start
    $seg1End   = $seg1.End

    # Segment 2 is appended right after segment 1.
    $seg1.Collapse(0)
    $seg1.InsertAfter("plus opłat")
    $seg2Start = $seg1End
    $seg2End   = $seg1.End

    # Segment 3 is appended right after segment 2.
    $seg1.Collapse(0)
    $seg1.InsertAfter("a")
    $seg3Start = $seg2End
    $seg3End   = $seg1.End

    # Segment 4 is appended right after segment 3.
    $seg1.Collapse(0)
    $seg1.InsertAfter(" oraz prowizja")
    $seg4Start = $seg3End
    $seg4End   = $seg1.End

    # Newly-inserted text automatically inherits the formatting of its
    # insertion point, so segments 2-4 currently all belong to one run
    # that is merged with segment 1. Force each of them into its own
    # <w:r> by nudging a formatting property away from, and back to,
    # its current value - this breaks the run apart without altering
    # the visible formatting (still Arial / bold / 12pt afterwards).
    $seg2 = $d.Range($seg2Start, $seg2End)
    $seg2.Bold = $false
    $seg2.Bold = $true

    $seg3 = $d.Range($seg3Start, $seg3End)
    $seg3.Bold = $false
    $seg3.Bold = $true

    $seg4 = $d.Range($seg4Start, $seg4End)
    $seg4.Bold = $false
    $seg4.Bold = $true
}

# ---------------------------------------------------------------------------
# 2) styles.xml: mark the built-in "Default Paragraph Font" character
#    style as semi-hidden (best effort - not all hosts expose a writable
#    property for w:semiHidden on a Style object).
# ---------------------------------------------------------------------------
try {
    $defStyle = $d.Styles("Domylnaczcionkaakapitu")
    $defStyle.Hidden = $true
} catch {
}
